$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.215.38"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "3.168.71"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.166.75"
$ws.Range("E8").Value = "  -1.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.548"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.83%  "
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("E11").Value = "  -9.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.519"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.05%  "
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").Value = "3.689.03"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "66.206.44"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "3.170.86"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "510.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.728"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.97%  "
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "501.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("E38").Value = "  -2.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0420"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.128"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.00%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0683"
$ws.Range("E42").Value = "  +6.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.296"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("D46").Value = "2.825.41"
$ws.Range("E46").Value = "  -3.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.21%  "
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.83%  "
